$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 3, pushing the existing "TestTwo" entry down to row 4
$ws.Rows("3:3").Insert()

# Set the new cell's value (becomes a new shared string entry)
$ws.Range("A3").Value = "TestSheetTwoFromMaster"

# Update the active selection to match the target workbook state
$ws.Range("E10").Select()
